$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Remove the existing hyperlink on the soon-to-move F2 cell and its style
$ws1.Hyperlinks.Delete()
$ws1.Range("F2").Style = "Normal"

# Copy header formatting (bold/border/center/top style) from the summary sheet header
$ws2.Range("A1:C1").Copy()
$ws1.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 2 (new notice): 第九届董事会第四十八次会议决议公告
$ws1.Range("A2").Value = "山东路桥"
$ws1.Range("B2").Value = "'000498"
$ws1.Range("C2").Value = "山东路桥:第九届董事会第四十八次会议决议公告"
$ws1.Range("D2").Value = "2023-03-24 00:00:00"
$ws1.Range("E2").Value = "2023-03-23 21:37:23:000"
$ws1.Range("F2").Value = "https://data.eastmoney.com/notices/detail/000498/AN202303231584517209.html"

# Row 3 (new notice): 关联交易公告
$ws1.Range("A3").Value = "山东路桥"
$ws1.Range("B3").Value = "'000498"
$ws1.Range("C3").Value = "山东路桥:关于中标荣乌高速烟威改扩建项目并投资灵犀七号及济南弘嘉的关联交易公告"
$ws1.Range("D3").Value = "2023-03-24 00:00:00"
$ws1.Range("E3").Value = "2023-03-23 21:37:07:000"
$ws1.Range("F3").Value = "https://data.eastmoney.com/notices/detail/000498/AN202303231584517208.html"

# Row 4 (new notice): 独立董事事前认可及独立意见
$ws1.Range("A4").Value = "山东路桥"
$ws1.Range("B4").Value = "'000498"
$ws1.Range("C4").Value = "山东路桥:独立董事关于第九届董事会第四十八次会议相关事项的事前认可及独立意见"
$ws1.Range("D4").Value = "2023-03-24 00:00:00"
$ws1.Range("E4").Value = "2023-03-23 21:37:02:000"
$ws1.Range("F4").Value = "https://data.eastmoney.com/notices/detail/000498/AN202303231584517221.html"

# Row 5 (original notice, shifted down from row 2): 可转换公司债券发行提示性公告
$ws1.Range("A5").Value = "山东路桥"
$ws1.Range("B5").Value = "'000498"
$ws1.Range("C5").Value = "山东路桥:向不特定对象发行可转换公司债券发行提示性公告"
$ws1.Range("D5").Value = "2023-03-24 00:00:00"
$ws1.Range("E5").Value = "2023-03-23 17:15:52:000"
$ws1.Range("F5").Value = "https://data.eastmoney.com/notices/detail/000498/AN202303231584507881.html"

# Clear the leading quote-prefix styling introduced by the apostrophe-escaped B values
$ws1.Range("B2:B5").Style = "Normal"

# Page margins: switch to Excel default (matches summary sheet)
$ws1.PageSetup.LeftMargin = 54
$ws1.PageSetup.RightMargin = 54
$ws1.PageSetup.TopMargin = 72
$ws1.PageSetup.BottomMargin = 72
$ws1.PageSetup.HeaderMargin = 36
$ws1.PageSetup.FooterMargin = 36

# Update the notice count on the summary sheet
$ws2.Range("C2").Value = 4

